# feat: add 2022-Q3 data
#
# - Duplicate the existing "2022-Q2" sheet to a new sheet right after it
#   (the duplicate keeps the old Q2 content/formatting and is renamed
#   back to "2022-Q2"), then replace the original sheet's content with
#   the new Q3 fund-holdings table and rename it to "2022-Q3". This
#   mirrors the upstream edit, where sheetId=2/rId2 becomes "2022-Q3"
#   and a fresh sheetId=3/rId3 becomes "2022-Q2".
# - Update the "总计" (summary) sheet: the new first data row holds the
#   Q3 totals, and the previous Q2 totals row shifts down one row.

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)
$oldQ2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1) Duplicate "2022-Q2" right after itself, then swap names so the
#    original sheet object (sheetId=2) becomes "2022-Q3" and the new
#    copy (sheetId=3) keeps the "2022-Q2" name/content.
# ---------------------------------------------------------------------
$oldQ2.Copy($null, $oldQ2)
$newQ2 = $wb.Worksheets.Item(3)

$oldQ2.Name = "2022-Q3"
$newQ2.Name = "2022-Q2"

$q3 = $oldQ2
$q3.Cells.Clear()

# A cell that already carries the shared header/index style (s="2") used
# throughout "总计" - copy it wherever that same style is needed so we
# reuse the existing style entry instead of inventing a new one.
$styleSrc = $summary.Range("B1")

# ---------------------------------------------------------------------
# 2) Populate "2022-Q3" with the fund holdings table.
# ---------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $styleSrc.Copy($q3.Cells.Item(1, $col))
    $q3.Cells.Item(1, $col).Value = $headers[$col - 2]
}

$rows = @(
    @("011460", "鹏华创新成长混合A",                               "10.08", "87.81", "3.75", "0.3780", 6),
    @("007484", "信澳核心科技混合",                                 "23.33", "94.48", "1.53", "0.3569", 8),
    @("501201", "红土创新科技创新 3 年封闭运作灵活配置混合",          "3.85",  "98.34", "4.38", "0.1686", 3),
    @("200010", "长城双动力混合A",                                 "3.29",  "93.10", "2.85", "0.0938", 9),
    @("014840", "招商裕华混合",                                    "2.04",  "63.58", "3.94", "0.0804", 6),
    @("015561", "长城双动力混合C",                                 "2.72",  "93.10", "2.85", "0.0775", 9),
    @("013721", "信澳景气优选混合A",                                "1.06",  "82.64", "6.90", "0.0731", 2),
    @("168401", "红土创新转型精选灵活配置混合（LOF）",               "0.87",  "92.08", "4.36", "0.0379", 1),
    @("011367", "创金合信群力一年定期开放混合（MOM）A",              "1.90",  "65.17", "1.89", "0.0359", 2),
    @("013722", "信澳景气优选混合C",                                "0.47",  "82.64", "6.90", "0.0324", 2),
    @("011461", "鹏华创新成长混合C",                                "0.37",  "87.81", "3.75", "0.0139", 6),
    @("011368", "创金合信群力一年定期开放混合（MOM）C",              "0.28",  "65.17", "1.89", "0.0053", 2)
)

$r = 2
foreach ($fund in $rows) {
    # Index column (A) shares the same style as the header row.
    $styleSrc.Copy($q3.Cells.Item($r, 1))
    $q3.Cells.Item($r, 1).Value = $r - 2

    # Fund code keeps any leading zeros -> force text.
    $q3.Cells.Item($r, 2).Value = "'" + $fund[0]
    $q3.Cells.Item($r, 2).Style = "Normal"

    # Fund name is plain text already.
    $q3.Cells.Item($r, 3).Value = $fund[1]

    # Scale / position / ratio / holding value are numeric-looking text
    # that must retain trailing zeros, so force text for these too.
    $q3.Cells.Item($r, 4).Value = "'" + $fund[2]
    $q3.Cells.Item($r, 4).Style = "Normal"
    $q3.Cells.Item($r, 5).Value = "'" + $fund[3]
    $q3.Cells.Item($r, 5).Style = "Normal"
    $q3.Cells.Item($r, 6).Value = "'" + $fund[4]
    $q3.Cells.Item($r, 6).Style = "Normal"
    $q3.Cells.Item($r, 7).Value = "'" + $fund[5]
    $q3.Cells.Item($r, 7).Style = "Normal"

    # Position rank is a real number.
    $q3.Cells.Item($r, 8).Value = $fund[6]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3) Update "总计": insert the Q3 totals as the new row 2, and push the
#    existing Q2 totals row down to row 3.
# ---------------------------------------------------------------------
# The current row 2 holds the Q2 totals ("2022-Q2", 1, 0.12); those move
# down to row 3 unchanged.

# Give row 3's index cell (A3) the same style as A2 (s="2"), then set its
# real value afterwards (Copy() would otherwise also copy A2's value).
$summary.Range("A2").Copy($summary.Range("A3"))
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.12

# Row 2 becomes the Q3 totals (A2 already holds 0 / the right style).
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 12
$summary.Range("D2").Value = 1.35

# Keep "总计" as the active sheet/tab, matching the original workbook.
$summary.Activate()
